$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the empty row 14, which shifts the old row 15
# ("Do we need more data?" note in column C) up to row 14.
$ws.Rows("14").Delete()

# Move the "Do we need more data?" note from column C to the new
# "Notes/Comments" column G on row 14, and clear the old cell.
$ws.Range("G14").Value = "Do we need more data?"
$ws.Range("C14").Value = ""

# Add the new "Notes/Comments" column header.
$ws.Range("G1").Value = "Notes/Comments"

# Annotate the row where continent data was missing (replaced with 'tbd').
$ws.Range("G5").Value = "no continent in orginal data, just 'tbd'"

# Annotate the row where year data was missing (replaced with 'n/a').
$ws.Range("G11").Value = "no year in original data, just 'n/a'"
